$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3132989.2
$ws.Range("J17").Value = 3132989.2
$ws.Range("L17").Value = 9398967.600000001
$ws.Range("N17").Value = -9399303.600000001
$ws.Range("H112").Value = 1154.1538
$ws.Range("I112").Value = 900
$ws.Range("J112").Value = 1160.8422
$ws.Range("K112").Value = 2700
$ws.Range("L112").Value = 3482.5266
$ws.Range("M112").Value = -1592
$ws.Range("N112").Value = -5698.5266
$ws.Range("H113").Value = 5828.077
$ws.Range("I113").Value = 2349.5417
$ws.Range("J113").Value = 11393.733
$ws.Range("K113").Value = 2349.5417
$ws.Range("L113").Value = 11393.733
$ws.Range("M113").Value = 904.4582999999998
$ws.Range("N113").Value = -17901.733
$ws.Range("H116").Value = 2272.889
$ws.Range("I116").Value = 2138.75
$ws.Range("J116").Value = 2380.2
$ws.Range("K116").Value = 2138.75
$ws.Range("L116").Value = 2380.2
$ws.Range("M116").Value = 1303.25
$ws.Range("N116").Value = -9264.200000000001
$ws.Range("H137").Value = 887.675
$ws.Range("I137").Value = 900.25
$ws.Range("J137").Value = 858.3333
$ws.Range("K137").Value = 2700.75
$ws.Range("L137").Value = 2574.9999
$ws.Range("M137").Value = -150.75
$ws.Range("N137").Value = -7674.9999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H122").Value = 1999.1428
$ws.Range("I122").Value = 1724.6666
$ws.Range("K122").Value = 5173.9998
$ws.Range("M122").Value = -2723.9998
$ws.Range("H131").Value = 39888
$ws.Range("J131").Value = 39888
$ws.Range("L131").Value = 39888
$ws.Range("N131").Value = -49968
$ws.Range("H132").Value = 2372
$ws.Range("I132").Value = 1601.4
$ws.Range("J132").Value = 3142.6
$ws.Range("K132").Value = 4804.200000000001
$ws.Range("L132").Value = 9427.799999999999
$ws.Range("M132").Value = -2274.200000000001
$ws.Range("N132").Value = -14487.8
$ws.Range("H139").Value = 40688.332
$ws.Range("J139").Value = 40688.332
$ws.Range("L139").Value = 40688.332
$ws.Range("N139").Value = -50968.332
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 373.2903
$ws.Range("I94").Value = 355.3125
$ws.Range("J94").Value = 392.46667
$ws.Range("K94").Value = 355.3125
$ws.Range("L94").Value = 392.46667
$ws.Range("M94").Value = 95.6875
$ws.Range("N94").Value = -1294.46667
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("N126").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3494.4736
$ws.Range("I31").Value = 1350.2456
$ws.Range("J31").Value = 9927.157999999999
$ws.Range("K31").Value = 1350.2456
$ws.Range("L31").Value = 9927.157999999999
$ws.Range("M31").Value = -1055.2456
$ws.Range("N31").Value = -10517.158
$ws.Range("H34").Value = 3494.4736
$ws.Range("I34").Value = 1350.2456
$ws.Range("J34").Value = 9927.157999999999
$ws.Range("K34").Value = 1350.2456
$ws.Range("L34").Value = 9927.157999999999
$ws.Range("M34").Value = -1148.2456
$ws.Range("N34").Value = -10331.158
$ws.Range("H99").Value = 16366
$ws.Range("I99").Value = 1912.4
$ws.Range("J99").Value = 52500
$ws.Range("K99").Value = 1912.4
$ws.Range("L99").Value = 52500
$ws.Range("M99").Value = -414.4000000000001
$ws.Range("N99").Value = -55496
$ws.Range("H126").Value = 16366
$ws.Range("I126").Value = 1912.4
$ws.Range("J126").Value = 52500
$ws.Range("K126").Value = 5737.200000000001
$ws.Range("L126").Value = 157500
$ws.Range("M126").Value = -3267.200000000001
$ws.Range("N126").Value = -162440
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 717053.2
$ws.Range("J107").Value = 908030.75
$ws.Range("L107").Value = 2724092.25
$ws.Range("N107").Value = -2727932.25
$ws.Range("H113").Value = 627.8333
$ws.Range("I113").Value = 572.5294
$ws.Range("J113").Value = 665.4400000000001
$ws.Range("K113").Value = 1717.5882
$ws.Range("L113").Value = 1996.32
$ws.Range("M113").Value = 452.4117999999999
$ws.Range("N113").Value = -6336.32
$ws.Range("H116").Value = 1375
$ws.Range("I116").Value = 823
$ws.Range("J116").Value = 2699.8
$ws.Range("K116").Value = 2469
$ws.Range("L116").Value = 8099.400000000001
$ws.Range("M116").Value = 973
$ws.Range("N116").Value = -14983.4
$ws.Range("H129").Value = 1431
$ws.Range("I129").Value = 791.4286
$ws.Range("J129").Value = 1804.0834
$ws.Range("K129").Value = 2374.2858
$ws.Range("L129").Value = 5412.2502
$ws.Range("M129").Value = 2625.7142
$ws.Range("N129").Value = -15412.2502
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5131464.5
$ws.Range("I102").Value = 7695295
$ws.Range("J102").Value = 3802.8
$ws.Range("K102").Value = 7695295
$ws.Range("L102").Value = 3802.8
$ws.Range("M102").Value = -7693673
$ws.Range("N102").Value = -7046.8
$ws.Range("H122").Value = 2667.8
$ws.Range("I122").Value = 2137.9583
$ws.Range("J122").Value = 3823.818
$ws.Range("K122").Value = 6413.874899999999
$ws.Range("L122").Value = 11471.454
$ws.Range("M122").Value = -3963.874899999999
$ws.Range("N122").Value = -16371.454
$ws.Range("H126").Value = 2091.0488
$ws.Range("I126").Value = 1877.96
$ws.Range("J126").Value = 2424
$ws.Range("K126").Value = 5633.88
$ws.Range("L126").Value = 7272
$ws.Range("M126").Value = -3163.88
$ws.Range("N126").Value = -12212
$ws.Range("H132").Value = 2716.2273
$ws.Range("I132").Value = 2504.923
$ws.Range("J132").Value = 3021.4443
$ws.Range("K132").Value = 7514.768999999999
$ws.Range("L132").Value = 9064.332900000001
$ws.Range("M132").Value = -4984.768999999999
$ws.Range("N132").Value = -14124.3329
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 46471.22
$ws.Range("I7").Value = 69152.60000000001
$ws.Range("J7").Value = 3943.625
$ws.Range("K7").Value = 69152.60000000001
$ws.Range("L7").Value = 3943.625
$ws.Range("M7").Value = -69040.60000000001
$ws.Range("N7").Value = -4167.625
$ws.Range("H40").Value = 44858.332
$ws.Range("I40").Value = 85583.336
$ws.Range("J40").Value = 4133.3335
$ws.Range("K40").Value = 85583.336
$ws.Range("L40").Value = 4133.3335
$ws.Range("M40").Value = -85447.336
$ws.Range("N40").Value = -4405.3335
$ws.Range("H122").Value = 2521.1538
$ws.Range("I122").Value = 1639
$ws.Range("J122").Value = 3072.5
$ws.Range("K122").Value = 4917
$ws.Range("L122").Value = 9217.5
$ws.Range("M122").Value = -2467
$ws.Range("N122").Value = -14117.5
$ws.Range("H126").Value = 46471.22
$ws.Range("I126").Value = 69152.60000000001
$ws.Range("J126").Value = 3943.625
$ws.Range("K126").Value = 207457.8
$ws.Range("L126").Value = 11830.875
$ws.Range("M126").Value = -204987.8
$ws.Range("N126").Value = -16770.875
$ws.Range("H132").Value = 9060.396000000001
$ws.Range("I132").Value = 7135.225
$ws.Range("J132").Value = 14984
$ws.Range("K132").Value = 21405.675
$ws.Range("L132").Value = 44952
$ws.Range("M132").Value = -18875.675
$ws.Range("N132").Value = -50012
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 16203.5
$ws.Range("I56").Value = 14250
$ws.Range("J56").Value = 18157
$ws.Range("K56").Value = 14250
$ws.Range("L56").Value = 18157
$ws.Range("M56").Value = -13536
$ws.Range("N56").Value = -19585
$ws.Range("H81").Value = 2712.5
$ws.Range("I81").Value = 2850
$ws.Range("J81").Value = 2692.8572
$ws.Range("K81").Value = 5700
$ws.Range("L81").Value = 5385.7144
$ws.Range("M81").Value = -4639
$ws.Range("N81").Value = -7507.7144
$ws.Range("H84").Value = 2712.5
$ws.Range("I84").Value = 2850
$ws.Range("J84").Value = 2692.8572
$ws.Range("K84").Value = 28500
$ws.Range("L84").Value = 26928.572
$ws.Range("M84").Value = -23196
$ws.Range("N84").Value = -37536.572
$ws.Range("H113").Value = 326.83334
$ws.Range("J113").Value = 334.33334
$ws.Range("L113").Value = 1003.00002
$ws.Range("N113").Value = -5343.00002
$ws.Range("H126").Value = 64184.375
$ws.Range("I126").Value = 143800
$ws.Range("J126").Value = 2261.111
$ws.Range("K126").Value = 431400
$ws.Range("L126").Value = 6783.333
$ws.Range("M126").Value = -428930
$ws.Range("N126").Value = -11723.333
